# Weekly fruit/vegetable price update: a new weekly record is inserted
# ahead of the existing "Ciboulette" history, pushing the old rows down
# by one and re-dating/re-pricing the new leading row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 174, shifting rows 174:181 down to 175:182
# (this also grows the used range / dimension to A1:R182).
$ws.Rows.Item(174).Insert()

# Populate the newly inserted row 174 with this week's record.
$ws.Cells.Item(174, 1).Value = 4
$ws.Cells.Item(174, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(174, 3).Value = "Los Lagos"
$ws.Cells.Item(174, 4).Value = 44610
$ws.Cells.Item(174, 5).Value = 10
$ws.Cells.Item(174, 6).Value = 100112039
$ws.Cells.Item(174, 7).Value = "Ciboulette"
$ws.Cells.Item(174, 8).Value = "Sin especificar"
$ws.Cells.Item(174, 9).Value = "Primera"
$ws.Cells.Item(174, 10).Value = 240
$ws.Cells.Item(174, 11).Value = 2500
$ws.Cells.Item(174, 12).Value = 3000
$ws.Cells.Item(174, 13).Value = 2750
$ws.Cells.Item(174, 14).Value = "$/docena de atados"
$ws.Cells.Item(174, 15).Value = "Región Metropolitana"
$ws.Cells.Item(174, 16).Value = 917
$ws.Cells.Item(174, 17).Value = 3
$ws.Cells.Item(174, 18).Value = "Hortaliza"
